# "upgrade left table until javakheti"
# Add the 2023 column (K) to the Mestia average-monthly-remuneration table,
# mirroring the formatting already used by the 2022 column (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header year + the three data rows (Total / Women / Men)
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1381.7
$ws.Range("K5").Value = 925.5
$ws.Range("K6").Value = 1629.1

# Copy column J's cell formatting (number format, font, fill, borders,
# alignment) onto the new column K cells without touching the values
# we just set.
$ws.Range("J3:J6").Copy()
$ws.Range("K3:K6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
